$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the text of the last existing entry (D25) - append a sentence about
#    the generalization of waitForEvent into waitForEventsTillTime.
$ws.Range("D25").Value = "Fix: Bad specification of ALL events - now timer events are still an OR condition. Implementation of waitForEventsTillTime by generalization of waitForEvent"

# 2. The effort logged on 2012-10-15 (row 25) grew from 1h to 2.5h.
$ws.Range("B25").Value = 2.5

# 3. Append a new day of work: 2012-10-16 (serial date 41198), 2.5h effort and
#    a new description string. Copy the date formatting from the cell above so
#    the new date cell keeps the same "ddd dd/mm/yyyy" display format.
$ws.Range("A26").Value = 41198
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B26").Value = 2.5
$ws.Range("D26").Value = "Code cleanup, suspendTillTime discarded. Implementation of enter/leaveCriticalSection"

# 4. Scroll the view down a bit and move the selection to the new first empty
#    row below the data (A27), matching the author's view position when the
#    sheet was saved.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A27").Select()
